$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 19: M2 Heat-Set Knurled Threaded Inserts - fill in MFR, MFR P/N, QTY
$ws.Range("C19").Value = "-"
$ws.Range("D19").Value = "-"
$ws.Range("E19").Value = 1

# Row 20: M2 Standoffs - fill in MFR, MFR P/N, QTY
$ws.Range("C20").Value = "-"
$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = 1

# Update the view: scroll back to top-left A1 and change selection to E21
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E21").Select()
